$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '62.961.00'
$ws.Range('E2').Value = '  +2.82%  '
$ws.Range('D3').Value = '2.958.92'
$ws.Range('E3').Value = '  +0.98%  '
$ws.Range('E4').Value = '  -0.02%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '595.13'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +0.26%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '146.52'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +1.19%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '1.00'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  +0.06%  '
$ws.Range('D8').Value = '2.958.64'
$ws.Range('E8').Value = '  +1.11%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.509'
$ws.Range('D9').Style = 'Normal'
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '7.27'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +3.86%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.152'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +6.40%  '
$ws.Range('E12').Value = '  +0.95%  '
$ws.Range('E13').Value = '  +6.30%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '33.13'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -1.29%  '
$ws.Range('E15').Value = '  -0.43%  '
$ws.Range('D16').Value = '3.448.19'
$ws.Range('E16').Value = '  +1.01%  '
$ws.Range('D17').Value = '62.819.95'
$ws.Range('E17').Value = '  +2.70%  '
$ws.Range('E18').Value = '  +0.07%  '
$ws.Range('D19').Value = '2.956.27'
$ws.Range('E19').Value = '  +1.06%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '443.41'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +2.41%  '
$ws.Range('E21').Value = '  +0.38%  '
$ws.Range('E22').Value = '  -1.57%  '
$ws.Range('E23').Value = '  -0.40%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '81.47'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -0.70%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '11.14'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +0.92%  '
$ws.Range('E26').Value = '  +0.61%  '
$ws.Range('E27').Value = '  -3.71%  '
$ws.Range('E28').Value = '  +0.05%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '7.32'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +4.89%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '2.63'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +0.90%  '
$ws.Range('E31').Value = '  -2.87%  '
$ws.Range('D32').Value = '0.0₃0978'
$ws.Range('E32').Value = '  +9.93%  '
$ws.Range('E33').Value = '  -1.13%  '
$ws.Range('E34').Value = '  -0.96%  '
$ws.Range('E35').Value = '  +0.02%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.995'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -1.67%  '
$ws.Range('E37').Value = '  +0.17%  '
$ws.Range('E38').Value = '  +3.84%  '
$ws.Range('E39').Value = '  +2.25%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '49.58'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -0.92%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '8.52'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -1.25%  '
$ws.Range('E42').Value = '  -4.90%  '
$ws.Range('E43').Value = '  -0.61%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '40.15'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -5.48%  '
$ws.Range('D45').Value = '2.725.39'
$ws.Range('E45').Value = '  +0.97%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '134.27'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +0.56%  '
$ws.Range('E47').Value = '  -2.83%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '361.90'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -2.30%  '
$ws.Range('E50').Value = '  -0.50%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '22.88'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -4.03%  '
